$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 19 (shifts old rows 19-30 down to 20-31,
# carrying their values/formatting with them, same as Excel's native
# "Insert Row" behaviour).
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the new weekly record.
$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(19, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(19, 4).Value = 44789
$ws.Cells.Item(19, 5).Value = 15
$ws.Cells.Item(19, 6).Value = 100112044
$ws.Cells.Item(19, 7).Value = "Perejil"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 300
$ws.Cells.Item(19, 11).Value = 1400
$ws.Cells.Item(19, 12).Value = 1500
$ws.Cells.Item(19, 13).Value = 1450
$ws.Cells.Item(19, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(19, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(19, 16).Value = 725
$ws.Cells.Item(19, 17).Value = 2
$ws.Cells.Item(19, 18).Value = "Hortaliza"
